$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels in A1 and B1: column A now holds the
# "TUR_GBR_841810_UV" series, column B now holds the "Date" series.
$hdrA = $ws.Range("A1").Value()
$hdrB = $ws.Range("B1").Value()
$ws.Range("A1").Value = $hdrB
$ws.Range("B1").Value = $hdrA

# For every data row, the date (previously in column A) moves to column B,
# and whatever value used to live in column B (if any) moves to column A.
for ($r = 2; $r -le 173; $r++) {
    $oldA = $ws.Cells.Item($r, 1).Value()
    $oldB = $ws.Cells.Item($r, 2).Value()

    if ($oldB -eq $null) {
        $ws.Cells.Item($r, 1).Value = $null
    } else {
        $ws.Cells.Item($r, 1).Value = $oldB
    }
    $ws.Cells.Item($r, 2).Value = $oldA
}
